$d = $word.ActiveDocument

# Anchor 1: the end of the "MÁQUINASEscola PRO-TEC" paragraph - everything
# after it (its own paragraph mark excluded) up through the end of the
# footer paragraph must be removed: a blank paragraph, the
# "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, and the
# "© 2020 ... Creative Commons Attribution" paragraph.
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("MÁQUINASEscola PRO-TEC", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$anchor.Collapse(0)        # collapse to just after the matched text
$anchor.MoveEnd(1, 1) | Out-Null   # extend over this paragraph's own mark

# Anchor 2: the end of the footer paragraph that contains the copyright /
# "Creative Commons Attribution" text - the last paragraph to be removed.
$tail = $d.Content.Duplicate
$tail.Find.Execute("Creative Commons Attribution", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$tail.Collapse(0)          # collapse to just after the matched text
$tail.MoveEnd(1, 1) | Out-Null     # extend over that paragraph's own mark

# Delete everything between the two anchors: the blank paragraph, the
# "Ver no Jupiter ..." paragraph and the "© 2020 ..." paragraph, while
# leaving the "MÁQUINASEscola PRO-TEC" paragraph and the paragraph that
# follows the removed block untouched.
$deleteRange = $d.Range($anchor.End, $tail.End)
$deleteRange.Delete()
